$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 1 : consolidate "Testing" " " "custom" " " "properties" into one run ---
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleFull = $titleRange.Characters(1, $titleRange.Length)
$titleFull.Text = "Testing custom properties"

# --- Subtitle 2 : consolidate runs around the two existing line breaks ---
$subShape = $s.Shapes.Item(2)
$subRange = $subShape.TextFrame.TextRange

# First segment: "This" " " "is" " " "a" " " "subtitle" -> "This is a subtitle"
$seg1 = $subRange.Characters(1, 18)
$seg1.Text = "This is a subtitle"

# Second segment (after the two line breaks): "A." " " "M." -> "A. M."
$seg2 = $subRange.Characters(21, 5)
$seg2.Text = "A. M."
